$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header style from D1 to E1 so the new header cell matches existing header formatting
$ws.Range("D1").Copy() | Out-Null
$ws.Range("E1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Populate column E with values (header + data rows)
$ws.Range("E1").Value = 'E'
$ws.Range("E2").Value = '35 B2Op'
$ws.Range("E3").Value = '36 B3Pop'
$ws.Range("E4").Value = '37 B4Pop'
$ws.Range("E5").Value = '38 B5C'
$ws.Range("E6").Value = '1 D1V'
$ws.Range("E7").Value = '3 D3V'
$ws.Range("E8").Value = '4 D4V'
$ws.Range("E9").Value = '5 D5V'
$ws.Range("E10").Value = '40 D2V'
$ws.Range("E11").Value = '42 D4V'
$ws.Range("E12").Value = '6 FB1C'
$ws.Range("E13").Value = '7 FB2C'
$ws.Range("E14").Value = '8 FB3C'
$ws.Range("E15").Value = '9 FB4C'
$ws.Range("E16").Value = '44 FB1C'
$ws.Range("E17").Value = '45 FB2C'
$ws.Range("E18").Value = '46 FB3Op'
$ws.Range("E19").Value = '47 FB4Pop'
$ws.Range("E20").Value = '48 FB5Pop'
$ws.Range("E21").Value = '10 H1C'
$ws.Range("E22").Value = '11 H2C'
$ws.Range("E23").Value = '12 H3C'
$ws.Range("E24").Value = '13 H4C'
$ws.Range("E25").Value = '14 H5C'
$ws.Range("E26").Value = '49 H1C'
$ws.Range("E27").Value = '51 H3C'
$ws.Range("E28").Value = '52 H4C'
$ws.Range("E29").Value = '53 H5De'
$ws.Range("E30").Value = '54 HH1De'
$ws.Range("E31").Value = '58 HH5De'
$ws.Range("E32").Value = '59 SF1C'
$ws.Range("E33").Value = '60 SF2C'
$ws.Range("E34").Value = '61 SF3C'
$ws.Range("E35").Value = '63 SF5C'
$ws.Range("E36").Value = '64 SLOp'
$ws.Range("E37").Value = '65 SOC'
$ws.Range("E38").Value = '66 ST1Rü'
$ws.Range("E39").Value = '67 ST2Rü'
$ws.Range("E40").Value = '68 ST3Rü'
$ws.Range("E41").Value = '69 ST4Rü'
$ws.Range("E42").Value = '67 ST5Rü'
$ws.Range("E43").Value = '15 Z1C'
$ws.Range("E44").Value = '16 Z2C'
$ws.Range("E45").Value = '17 Z3C'
$ws.Range("E46").Value = '18 Z4C'
$ws.Range("E47").Value = '19 Z5C'
